# Apply updated cryptocurrency market data to sheet1 (cryptos.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be misread as a number (losing
# trailing zeros / exact decimal representation) are forced to Text format
# first, then restored to the default "Normal" style so no stray number
# format sticks around on the cell.

$ws.Range('D2').Value = '97.159.19'
$ws.Range('E2').Value = '  +0.56%  '
$ws.Range('D3').Value = '3.708.36'
$ws.Range('E3').Value = '  +0.78%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '237.65'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.50%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.91'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.15%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '660.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.60%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('B9').Value = 'USDC'
$ws.Range('C9').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.00%  '
$ws.Range('B10').Value = 'Cardano'
$ws.Range('C10').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.07'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.35%  '
$ws.Range('D11').Value = '3.707.79'
$ws.Range('E11').Value = '  +0.79%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000318'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +18.11%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '44.25'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.66%  '
$ws.Range('E14').Value = '  +1.71%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.80'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.42%  '
$ws.Range('D16').Value = '4.399.27'
$ws.Range('E16').Value = '  +0.81%  '
$ws.Range('D17').Value = '97.039.70'
$ws.Range('E17').Value = '  +0.76%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '9.13'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.66%  '
$ws.Range('D19').Value = '3.695.05'
$ws.Range('E19').Value = '  +0.43%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.03'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.12%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '18.71'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.17%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.504'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.01%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '522.44'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.14%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.43'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.34%  '
$ws.Range('E25').Value = '  +4.21%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.91'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.34%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '101.75'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.18%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.193'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +13.77%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '13.59'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.01%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '12.76'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.94%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.06'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.25%  '
$ws.Range('B32').Value = 'Dai'
$ws.Range('C32').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.00'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.04%  '
$ws.Range('B33').Value = 'Cronos'
$ws.Range('C33').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.190'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.56%  '
$ws.Range('B34').Value = 'Fetch.AI'
$ws.Range('C34').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.89'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.00%  '
$ws.Range('B35').Value = 'Binance-PegBSC-USD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('B36').Value = 'EthereumClassic'
$ws.Range('C36').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '32.31'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.60%  '
$ws.Range('B37').Value = 'Bittensor'
$ws.Range('C37').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '653.26'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.64%  '
$ws.Range('B38').Value = 'PolygonEcosystemToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.595'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.14%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.86'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.00%  '
$ws.Range('B40').Value = 'USDe'
$ws.Range('C40').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.00'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.02%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.165'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.17%  '
$ws.Range('B42').Value = 'ImmutableX'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.05'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.68%  '
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.80'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.93%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '40.53'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -9.23%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.484'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +7.06%  '
$ws.Range('B46').Value = 'ARBITRUM'
$ws.Range('C46').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.972'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.48%  '
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0460'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.52%  '
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.31'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.39%  '
$ws.Range('B49').Value = 'WhiteBITCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '23.61'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.18%  '
$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.72'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.05%  '
$ws.Range('B51').Value = 'MantraDAO'
$ws.Range('C51').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.53'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.44%  '
